# Add I/O scaling factors input variable (#163)
#
# This adds a new row to the "Key to Variables" sheet describing the
# "BPCiObIC" (BAU Percent Change in Output by ISIC Code) variable, inserted
# right above the existing "BPEaCP" row (old row 158), and updates the
# active sheet / selection to reflect where the author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new row at position 158 (pushes the old row 158.. down by one,
# inheriting the formatting of the row above it, same as Excel's UI
# "Insert Sheet Rows").
$ws.Rows.Item(158).Insert()

# Populate the new row with the new variable's data.
$ws.Range("A158").Value = "io-model"
$ws.Range("B158").Value = "BPCiObIC"
$ws.Range("C158").Value = "BAU Percent Change in Output by ISIC Code"
$ws.Range("E158").Value = "BVAbIC"
$ws.Range("F158").Value = "high"

# Make "Key to Variables" the active sheet/tab, and leave the selection on
# the newly added row, scrolled down near it.
$ws.Activate()
$ws.Range("A158").Select()
